$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the "MZ" (Mozambique) country row, being updated for Moz training.
# gridPriceUrb / gridPriceRur / devicePrice
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 30000

# cookStove1 / cookStovePrice1 - replace old rich-text "COUNTRY SPECIFIC ..." label
# with a plain cookstove name and new price.
$ws.Range("F2").Value = "Chitetezo wood burner"
$ws.Range("G2").Value = 300

# cookStove3 / cookStovePrice3 - same for the charcoal stove.
$ws.Range("H2").Value = "Envirofit CH-2200 charcool cookstove"
$ws.Range("I2").Value = 3000

# Leave the selection where the editor ended up.
[void]$ws.Range("C2").Select()
